$d = $word.ActiveDocument

# Phase 1: replace each original requisito line with a unique placeholder
# to avoid collisions since this is a full reordering (permutation) of the list.
$d.Content.Find.Execute("LOB1053 -  Física III  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER00@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER01@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1008 -  Ciência, Tecnologia e Sociedade  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER02@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1036 -  Geometria Analítica  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER03@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1037 -  Àlgebra Linear  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER04@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1041 -  Física Experimental II  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER05@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1042 -  Física Experimental IV  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER06@@", 2) | Out-Null
$d.Content.Find.Execute("LOQ4095 -  Química Geral Experimental  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER07@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1039 -  Física Experimental III  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER08@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1018 -  Física I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER09@@", 2) | Out-Null
$d.Content.Find.Execute("LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER10@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1004 -  Cálculo II  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER11@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1038 -  Física Experimental I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER12@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1052 -  Cálculo III  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER13@@", 2) | Out-Null
$d.Content.Find.Execute("LOM3236 -  Processos de Fabricação  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER14@@", 2) | Out-Null
$d.Content.Find.Execute("LOM3261 -  Métodos Numéricos e Aplicações  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER15@@", 2) | Out-Null
$d.Content.Find.Execute("LOM3218 -  Introdução à Engenharia Física  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER16@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1003 -  Cálculo I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER17@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1006 -  Cálculo IV  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER18@@", 2) | Out-Null
$d.Content.Find.Execute("LOM3241 -  Química de Materiais  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER19@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1021 -  Física IV  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER20@@", 2) | Out-Null
$d.Content.Find.Execute("LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER21@@", 2) | Out-Null
$d.Content.Find.Execute("LOM3260 -  Computação Científica em Python  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER22@@", 2) | Out-Null
$d.Content.Find.Execute("LOM3204 -  Desenho Técnico e Projeto Assistido por Computador  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER23@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1012 -  Estatística  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER24@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1019 -  Física II  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER25@@", 2) | Out-Null

# Phase 2: replace each placeholder with the final (reordered) text
$d.Content.Find.Execute("@@PLACEHOLDER00@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER01@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3261 -  Métodos Numéricos e Aplicações  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER02@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4095 -  Química Geral Experimental  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER03@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER04@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1006 -  Cálculo IV  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER05@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1008 -  Ciência, Tecnologia e Sociedade  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER06@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1037 -  Àlgebra Linear  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER07@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1053 -  Física III  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER08@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3241 -  Química de Materiais  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER09@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1003 -  Cálculo I  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER10@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1012 -  Estatística  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER11@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1018 -  Física I  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER12@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1036 -  Geometria Analítica  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER13@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1038 -  Física Experimental I  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER14@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1039 -  Física Experimental III  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER15@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1041 -  Física Experimental II  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER16@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1042 -  Física Experimental IV  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER17@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1052 -  Cálculo III  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER18@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER19@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3218 -  Introdução à Engenharia Física  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER20@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3260 -  Computação Científica em Python  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER21@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1004 -  Cálculo II  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER22@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1019 -  Física II  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER23@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1021 -  Física IV  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER24@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3204 -  Desenho Técnico e Projeto Assistido por Computador  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER25@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3236 -  Processos de Fabricação  (Requisito)", 2) | Out-Null

Write-Output "done"
